$d = $word.ActiveDocument

# --- Change 1: Title paragraph — append a bold " – Maciej Morgalla" run ---
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$insertPoint = $d.Range($titleRange.Start + $titleRange.Text.TrimEnd([char]13,[char]7).Length, $titleRange.Start + $titleRange.Text.TrimEnd([char]13,[char]7).Length)
$insertPoint.InsertAfter(" – Maciej Morgalla")
# Force a distinct run boundary (same bold formatting as the preceding run)
$insertPoint.Font.Bold = $false
$insertPoint.Font.Bold = $true

# --- Change 2: "Funkcja celu" bullet — replace "wygladu" with "dzialania" ---
$findRange = $d.Content
$findRange.Find.Execute("wyglądu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Text = "działania"
# Force the replaced word onto its own run (split away from the surrounding text)
$findRange.Font.Bold = $true
$findRange.Font.Bold = $false
